# Sort TM/HM moves automatically, work on Clone Charmander/Charmeleon
#
# Appends the remaining TM moves (as new TM_MOVE rows) to the bottom of the
# move tables on the Clone_Charmander and Clone_Charmeleon sheets.

$wb = $excel.ActiveWorkbook

$newTmMoves = @(
    "Dragon Claw",
    "Hidden Power",
    "Frustration",
    "Iron Tail",
    "Dig",
    "Brick Break",
    "Facade",
    "Secret Power",
    "Rest",
    "Cut",
    "Strength",
    "Rock Smash"
)

$sheetNames = @("Clone_Charmander", "Clone_Charmeleon")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Last existing TM_MOVE row is row 79 (A79 = "TM_MOVE", C79 = "Toxic").
    $startRow = 80
    $row = $startRow

    foreach ($moveName in $newTmMoves) {
        $ws.Cells.Item($row, 1).Value = "TM_MOVE"
        $ws.Cells.Item($row, 3).Value = $moveName
        $row = $row + 1
    }

    $lastRow = $row - 1

    [void]$ws.Activate()
    [void]$ws.Range("C" + $lastRow).Select()
}

# Put the selection/active cell on Clone_Charmeleon the same way the
# original edit left it (column A, near the newly appended rows), and make
# sure it ends up as the active sheet again.
$wsCharmeleon = $wb.Worksheets.Item("Clone_Charmeleon")
[void]$wsCharmeleon.Activate()
[void]$wsCharmeleon.Range("A76").Select()
